$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 917.2308
$ws.Cells.Item(2, 10).Value = 1492.7142
$ws.Cells.Item(2, 12).Value = 1492.7142
$ws.Cells.Item(2, 14).Value = -1718.7142
$ws.Cells.Item(19, 8).Value = 3479.8
$ws.Cells.Item(19, 9).Value = 2000
$ws.Cells.Item(19, 10).Value = 3849.75
$ws.Cells.Item(19, 11).Value = 2000
$ws.Cells.Item(19, 12).Value = 3849.75
$ws.Cells.Item(19, 13).Value = -1825
$ws.Cells.Item(19, 14).Value = -4199.75
$ws.Cells.Item(70, 8).Value = 8653.35
$ws.Cells.Item(70, 10).Value = 10963.692
$ws.Cells.Item(70, 12).Value = 32891.076
$ws.Cells.Item(70, 14).Value = -33431.076
$ws.Cells.Item(73, 8).Value = 8653.35
$ws.Cells.Item(73, 10).Value = 10963.692
$ws.Cells.Item(73, 12).Value = 32891.076
$ws.Cells.Item(73, 14).Value = -34763.076
$ws.Cells.Item(80, 8).Value = 3343.2778
$ws.Cells.Item(80, 9).Value = 3141
$ws.Cells.Item(80, 10).Value = 3505.1
$ws.Cells.Item(80, 11).Value = 9423
$ws.Cells.Item(80, 12).Value = 10515.3
$ws.Cells.Item(80, 13).Value = -8425
$ws.Cells.Item(80, 14).Value = -12511.3
$ws.Cells.Item(83, 8).Value = 3343.2778
$ws.Cells.Item(83, 9).Value = 3141
$ws.Cells.Item(83, 10).Value = 3505.1
$ws.Cells.Item(83, 11).Value = 28269
$ws.Cells.Item(83, 12).Value = 31545.9
$ws.Cells.Item(83, 13).Value = -23277
$ws.Cells.Item(83, 14).Value = -41529.89999999999
$ws.Cells.Item(86, 8).Value = 5621.9
$ws.Cells.Item(86, 9).Value = 3579.75
$ws.Cells.Item(86, 10).Value = 6983.3335
$ws.Cells.Item(86, 11).Value = 3579.75
$ws.Cells.Item(86, 12).Value = 6983.3335
$ws.Cells.Item(86, 13).Value = -2456.75
$ws.Cells.Item(86, 14).Value = -9229.333500000001
$ws.Cells.Item(89, 8).Value = 5621.9
$ws.Cells.Item(89, 9).Value = 3579.75
$ws.Cells.Item(89, 10).Value = 6983.3335
$ws.Cells.Item(89, 11).Value = 17898.75
$ws.Cells.Item(89, 12).Value = 34916.6675
$ws.Cells.Item(89, 13).Value = -12282.75
$ws.Cells.Item(89, 14).Value = -46148.6675
$ws.Cells.Item(98, 8).Value = 4723.409
$ws.Cells.Item(98, 9).Value = 2093.4666
$ws.Cells.Item(98, 10).Value = 10359
$ws.Cells.Item(98, 11).Value = 2093.4666
$ws.Cells.Item(98, 12).Value = 10359
$ws.Cells.Item(98, 13).Value = -595.4666000000002
$ws.Cells.Item(98, 14).Value = -13355
$ws.Cells.Item(113, 8).Value = 4913.3887
$ws.Cells.Item(113, 9).Value = 4372.2856
$ws.Cells.Item(113, 10).Value = 5670.933
$ws.Cells.Item(113, 11).Value = 4372.2856
$ws.Cells.Item(113, 12).Value = 5670.933
$ws.Cells.Item(113, 13).Value = -1118.2856
$ws.Cells.Item(113, 14).Value = -12178.933
$ws.Cells.Item(122, 8).Value = 4723.409
$ws.Cells.Item(122, 9).Value = 2093.4666
$ws.Cells.Item(122, 10).Value = 10359
$ws.Cells.Item(122, 11).Value = 6280.399800000001
$ws.Cells.Item(122, 12).Value = 31077
$ws.Cells.Item(122, 13).Value = -3830.399800000001
$ws.Cells.Item(122, 14).Value = -35977
$ws.Cells.Item(127, 8).Value = 4453.4443
$ws.Cells.Item(127, 9).Value = 4635.125
$ws.Cells.Item(127, 10).Value = 3000
$ws.Cells.Item(127, 11).Value = 13905.375
$ws.Cells.Item(127, 12).Value = 9000
$ws.Cells.Item(127, 13).Value = -8945.375
$ws.Cells.Item(127, 14).Value = -18920
$ws.Cells.Item(132, 8).Value = 14496732
$ws.Cells.Item(132, 9).Value = 20412108
$ws.Cells.Item(132, 11).Value = 61236324
$ws.Cells.Item(132, 13).Value = -61233794
$ws.Cells.Item(133, 8).Value = 87000
$ws.Cells.Item(133, 10).Value = 87000
$ws.Cells.Item(133, 12).Value = 87000
$ws.Cells.Item(133, 14).Value = -97120
$ws.Cells.Item(137, 8).Value = 1547.262
$ws.Cells.Item(137, 9).Value = 1046.3214
$ws.Cells.Item(137, 11).Value = 3138.9642
$ws.Cells.Item(137, 13).Value = -588.9642000000003
$ws.Cells.Item(138, 8).Value = 4913.7925
$ws.Cells.Item(138, 9).Value = 1364.0435
$ws.Cells.Item(138, 10).Value = 7635.2666
$ws.Cells.Item(138, 11).Value = 4092.1305
$ws.Cells.Item(138, 12).Value = 22905.7998
$ws.Cells.Item(138, 13).Value = 1047.8695
$ws.Cells.Item(138, 14).Value = -33185.7998
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3527.923
$ws.Cells.Item(45, 9).Value = 1985.6666
$ws.Cells.Item(45, 11).Value = 1985.6666
$ws.Cells.Item(45, 13).Value = -1608.6666
$ws.Cells.Item(61, 8).Value = 5717.375
$ws.Cells.Item(61, 9).Value = 2515.8333
$ws.Cells.Item(61, 11).Value = 2515.8333
$ws.Cells.Item(61, 13).Value = -2303.8333
$ws.Cells.Item(88, 8).Value = 2811.7727
$ws.Cells.Item(88, 9).Value = 1932
$ws.Cells.Item(88, 10).Value = 3141.6875
$ws.Cells.Item(88, 11).Value = 1932
$ws.Cells.Item(88, 12).Value = 3141.6875
$ws.Cells.Item(88, 13).Value = -1526
$ws.Cells.Item(88, 14).Value = -3953.6875
$ws.Cells.Item(91, 8).Value = 2811.7727
$ws.Cells.Item(91, 9).Value = 1932
$ws.Cells.Item(91, 10).Value = 3141.6875
$ws.Cells.Item(91, 11).Value = 1932
$ws.Cells.Item(91, 12).Value = 3141.6875
$ws.Cells.Item(91, 13).Value = -528
$ws.Cells.Item(91, 14).Value = -5949.6875
$ws.Cells.Item(122, 8).Value = 8853.833000000001
$ws.Cells.Item(122, 9).Value = 8853.833000000001
$ws.Cells.Item(122, 11).Value = 26561.499
$ws.Cells.Item(122, 13).Value = -24111.499
$ws.Cells.Item(132, 8).Value = 7958.727
$ws.Cells.Item(132, 10).Value = 10285.643
$ws.Cells.Item(132, 12).Value = 30856.929
$ws.Cells.Item(132, 14).Value = -35916.929
$ws.Cells.Item(136, 8).Value = 5717.375
$ws.Cells.Item(136, 9).Value = 2515.8333
$ws.Cells.Item(136, 11).Value = 7547.499899999999
$ws.Cells.Item(136, 13).Value = -4997.499899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1744.037
$ws.Cells.Item(20, 10).Value = 1764.95
$ws.Cells.Item(20, 12).Value = 1764.95
$ws.Cells.Item(20, 14).Value = -2258.95
$ws.Cells.Item(82, 8).Value = 12869.318
$ws.Cells.Item(82, 10).Value = 53672.668
$ws.Cells.Item(82, 12).Value = 53672.668
$ws.Cells.Item(82, 14).Value = -54438.668
$ws.Cells.Item(85, 8).Value = 12869.318
$ws.Cells.Item(85, 10).Value = 53672.668
$ws.Cells.Item(85, 12).Value = 53672.668
$ws.Cells.Item(85, 14).Value = -56324.668
$ws.Cells.Item(86, 8).Value = 8335444
$ws.Cells.Item(86, 9).Value = 22224220
$ws.Cells.Item(86, 10).Value = 2178.8
$ws.Cells.Item(86, 11).Value = 22224220
$ws.Cells.Item(86, 12).Value = 2178.8
$ws.Cells.Item(86, 13).Value = -22223097
$ws.Cells.Item(86, 14).Value = -4424.8
$ws.Cells.Item(89, 8).Value = 8335444
$ws.Cells.Item(89, 9).Value = 22224220
$ws.Cells.Item(89, 10).Value = 2178.8
$ws.Cells.Item(89, 11).Value = 111121100
$ws.Cells.Item(89, 12).Value = 10894
$ws.Cells.Item(89, 13).Value = -111115484
$ws.Cells.Item(89, 14).Value = -22126
$ws.Cells.Item(94, 8).Value = 1256.0555
$ws.Cells.Item(94, 9).Value = 715.3
$ws.Cells.Item(94, 11).Value = 715.3
$ws.Cells.Item(94, 13).Value = -264.3
$ws.Cells.Item(134, 8).Value = 1689.8182
$ws.Cells.Item(134, 9).Value = 1500.5
$ws.Cells.Item(134, 11).Value = 4501.5
$ws.Cells.Item(134, 13).Value = -1966.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3115.6924
$ws.Cells.Item(31, 9).Value = 2033.15
$ws.Cells.Item(31, 10).Value = 4255.2104
$ws.Cells.Item(31, 11).Value = 2033.15
$ws.Cells.Item(31, 12).Value = 4255.2104
$ws.Cells.Item(31, 13).Value = -1738.15
$ws.Cells.Item(31, 14).Value = -4845.2104
$ws.Cells.Item(32, 8).Value = 2211.75
$ws.Cells.Item(32, 9).Value = 2211.75
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 2211.75
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -1895.75
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 3115.6924
$ws.Cells.Item(34, 9).Value = 2033.15
$ws.Cells.Item(34, 10).Value = 4255.2104
$ws.Cells.Item(34, 11).Value = 2033.15
$ws.Cells.Item(34, 12).Value = 4255.2104
$ws.Cells.Item(34, 13).Value = -1831.15
$ws.Cells.Item(34, 14).Value = -4659.2104
$ws.Cells.Item(99, 8).Value = 3147.3333
$ws.Cells.Item(99, 9).Value = 3291.9
$ws.Cells.Item(99, 10).Value = 2424.5
$ws.Cells.Item(99, 11).Value = 3291.9
$ws.Cells.Item(99, 12).Value = 2424.5
$ws.Cells.Item(99, 13).Value = -1793.9
$ws.Cells.Item(99, 14).Value = -5420.5
$ws.Cells.Item(126, 8).Value = 3147.3333
$ws.Cells.Item(126, 9).Value = 3291.9
$ws.Cells.Item(126, 10).Value = 2424.5
$ws.Cells.Item(126, 11).Value = 9875.700000000001
$ws.Cells.Item(126, 12).Value = 7273.5
$ws.Cells.Item(126, 13).Value = -7405.700000000001
$ws.Cells.Item(126, 14).Value = -12213.5
$ws.Cells.Item(134, 8).Value = 1556.5217
$ws.Cells.Item(134, 9).Value = 1521.2
$ws.Cells.Item(134, 10).Value = 1792
$ws.Cells.Item(134, 11).Value = 4563.6
$ws.Cells.Item(134, 12).Value = 5376
$ws.Cells.Item(134, 13).Value = -2028.6
$ws.Cells.Item(134, 14).Value = -10446
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 9286.546
$ws.Cells.Item(38, 9).Value = 440.4
$ws.Cells.Item(38, 10).Value = 16658.334
$ws.Cells.Item(38, 11).Value = 1321.2
$ws.Cells.Item(38, 12).Value = 49975.00199999999
$ws.Cells.Item(38, 13).Value = -974.1999999999998
$ws.Cells.Item(38, 14).Value = -50669.00199999999
$ws.Cells.Item(75, 8).Value = 5907.778
$ws.Cells.Item(75, 9).Value = 4436
$ws.Cells.Item(75, 10).Value = 6202.1333
$ws.Cells.Item(75, 11).Value = 13308
$ws.Cells.Item(75, 12).Value = 18606.3999
$ws.Cells.Item(75, 13).Value = -12310
$ws.Cells.Item(75, 14).Value = -20602.3999
$ws.Cells.Item(78, 8).Value = 5907.778
$ws.Cells.Item(78, 9).Value = 4436
$ws.Cells.Item(78, 10).Value = 6202.1333
$ws.Cells.Item(78, 11).Value = 39924
$ws.Cells.Item(78, 12).Value = 55819.1997
$ws.Cells.Item(78, 13).Value = -34932
$ws.Cells.Item(78, 14).Value = -65803.1997
$ws.Cells.Item(88, 8).Value = 7197.5
$ws.Cells.Item(88, 10).Value = 7197.5
$ws.Cells.Item(88, 12).Value = 21592.5
$ws.Cells.Item(88, 14).Value = -22448.5
$ws.Cells.Item(91, 8).Value = 7197.5
$ws.Cells.Item(91, 10).Value = 7197.5
$ws.Cells.Item(91, 12).Value = 21592.5
$ws.Cells.Item(91, 14).Value = -24556.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 504.5
$ws.Cells.Item(2, 9).Value = 504.5
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 504.5
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -391.5
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 9742.571
$ws.Cells.Item(80, 9).Value = 15332.667
$ws.Cells.Item(80, 10).Value = 5550
$ws.Cells.Item(80, 11).Value = 15332.667
$ws.Cells.Item(80, 12).Value = 5550
$ws.Cells.Item(80, 13).Value = -14334.667
$ws.Cells.Item(80, 14).Value = -7546
$ws.Cells.Item(83, 8).Value = 9742.571
$ws.Cells.Item(83, 9).Value = 15332.667
$ws.Cells.Item(83, 10).Value = 5550
$ws.Cells.Item(83, 11).Value = 76663.33499999999
$ws.Cells.Item(83, 12).Value = 27750
$ws.Cells.Item(83, 13).Value = -71671.33499999999
$ws.Cells.Item(83, 14).Value = -37734
$ws.Cells.Item(107, 8).Value = 752.65216
$ws.Cells.Item(107, 9).Value = 503
$ws.Cells.Item(107, 11).Value = 503
$ws.Cells.Item(107, 13).Value = 1417
$ws.Cells.Item(122, 8).Value = 3948
$ws.Cells.Item(122, 9).Value = 3930.6667
$ws.Cells.Item(122, 10).Value = 4000
$ws.Cells.Item(122, 11).Value = 11792.0001
$ws.Cells.Item(122, 12).Value = 12000
$ws.Cells.Item(122, 13).Value = -9342.000100000001
$ws.Cells.Item(122, 14).Value = -16900
$ws.Cells.Item(126, 8).Value = 15575.65
$ws.Cells.Item(126, 9).Value = 5592.875
$ws.Cells.Item(126, 11).Value = 16778.625
$ws.Cells.Item(126, 13).Value = -14308.625
$ws.Cells.Item(132, 8).Value = 3678.7605
$ws.Cells.Item(132, 9).Value = 4043.5688
$ws.Cells.Item(132, 10).Value = 2051.1538
$ws.Cells.Item(132, 11).Value = 12130.7064
$ws.Cells.Item(132, 12).Value = 6153.4614
$ws.Cells.Item(132, 13).Value = -9600.706399999999
$ws.Cells.Item(132, 14).Value = -11213.4614
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2329.476
$ws.Cells.Item(16, 9).Value = 2021.5
$ws.Cells.Item(16, 10).Value = 2740.111
$ws.Cells.Item(16, 11).Value = 2021.5
$ws.Cells.Item(16, 12).Value = 2740.111
$ws.Cells.Item(16, 13).Value = -1851.5
$ws.Cells.Item(16, 14).Value = -3080.111
$ws.Cells.Item(55, 8).Value = 1470.7894
$ws.Cells.Item(55, 9).Value = 664
$ws.Cells.Item(55, 10).Value = 2057.5454
$ws.Cells.Item(55, 11).Value = 664
$ws.Cells.Item(55, 12).Value = 2057.5454
$ws.Cells.Item(55, 13).Value = -491
$ws.Cells.Item(55, 14).Value = -2403.5454
$ws.Cells.Item(61, 8).Value = 9890.538
$ws.Cells.Item(61, 10).Value = 7131.8
$ws.Cells.Item(61, 12).Value = 7131.8
$ws.Cells.Item(61, 14).Value = -7535.8
$ws.Cells.Item(82, 8).Value = 17687.666
$ws.Cells.Item(82, 10).Value = 1137.7142
$ws.Cells.Item(82, 12).Value = 1137.7142
$ws.Cells.Item(82, 14).Value = -1859.7142
$ws.Cells.Item(85, 8).Value = 17687.666
$ws.Cells.Item(85, 10).Value = 1137.7142
$ws.Cells.Item(85, 12).Value = 1137.7142
$ws.Cells.Item(85, 14).Value = -3633.7142
$ws.Cells.Item(109, 8).Value = 85000
$ws.Cells.Item(109, 10).Value = 85000
$ws.Cells.Item(109, 12).Value = 85000
$ws.Cells.Item(109, 14).Value = -87774
$ws.Cells.Item(113, 8).Value = 9890.538
$ws.Cells.Item(113, 10).Value = 7131.8
$ws.Cells.Item(113, 12).Value = 7131.8
$ws.Cells.Item(113, 14).Value = -11471.8
$ws.Cells.Item(122, 8).Value = 12779.45
$ws.Cells.Item(122, 9).Value = 13491.846
$ws.Cells.Item(122, 10).Value = 11456.429
$ws.Cells.Item(122, 11).Value = 40475.538
$ws.Cells.Item(122, 12).Value = 34369.287
$ws.Cells.Item(122, 13).Value = -38025.538
$ws.Cells.Item(122, 14).Value = -39269.287
$ws.Cells.Item(132, 8).Value = 5634.5454
$ws.Cells.Item(132, 9).Value = 4000
$ws.Cells.Item(132, 11).Value = 12000
$ws.Cells.Item(132, 13).Value = -9470
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4739.8423
$ws.Cells.Item(81, 10).Value = 9670.714
$ws.Cells.Item(81, 12).Value = 19341.428
$ws.Cells.Item(81, 14).Value = -21463.428
$ws.Cells.Item(84, 8).Value = 4739.8423
$ws.Cells.Item(84, 10).Value = 9670.714
$ws.Cells.Item(84, 12).Value = 96707.14
$ws.Cells.Item(84, 14).Value = -107315.14
$ws.Cells.Item(108, 8).Value = 100626
$ws.Cells.Item(108, 10).Value = 100626
$ws.Cells.Item(108, 12).Value = 100626
$ws.Cells.Item(108, 14).Value = -108306
$ws.Cells.Item(132, 8).Value = 1777.6072
$ws.Cells.Item(132, 9).Value = 1685.3636
$ws.Cells.Item(132, 10).Value = 2115.8333
$ws.Cells.Item(132, 11).Value = 5056.0908
$ws.Cells.Item(132, 12).Value = 6347.499899999999
$ws.Cells.Item(132, 13).Value = -2526.0908
$ws.Cells.Item(132, 14).Value = -11407.4999
$ws.Cells.Item(136, 8).Value = 3018.029
$ws.Cells.Item(136, 9).Value = 1572.3823
$ws.Cells.Item(136, 11).Value = 4717.1469
$ws.Cells.Item(136, 13).Value = -2167.1469
